# MasterSaveDemo / Database_original.xlsx edit
# Commit: "Fix Trung bug: BaoCaoMoDong_ViewModel, Tien: ThayDoiQuyDinh_ViewModel Set input (only number)"
#
# Summary of the change:
#  - LOAITIETKIEM (sheet "LOAITIETKIEM"): columns F/G ("rut het?" / "co hieu luc?")
#    switch from raw 0/1 numbers to descriptive text values, row 1 gets taller,
#    column F gets a wider custom width, and the selection moves to F2.
#  - SOTIETKIEM (sheet "SOTIETKIEM"): column F (so CMND/so dien thoai ...) is
#    reformatted as Text (@) so big numbers are no longer mangled, one value is
#    corrected, and the sheet becomes the active tab/view.
#  - THAMSO loses the "active tab" marker (SOTIETKIEM becomes active instead).
#  - PHANQUYEN selection resets back to the top of its used range.

$wb = $excel.ActiveWorkbook

$wsLoaiTietKiem = $wb.Worksheets.Item("LOAITIETKIEM")
$wsSoTietKiem   = $wb.Worksheets.Item("SOTIETKIEM")
$wsThamSo       = $wb.Worksheets.Item("THAMSO")
$wsPhanQuyen    = $wb.Worksheets.Item("PHANQUYEN")

# ---------------------------------------------------------------------------
# LOAITIETKIEM: F/G columns become descriptive text instead of 0/1 flags.
# ---------------------------------------------------------------------------

# NOTE: shared-string indices are assigned in first-seen order, so "Có" must
# be written before "Rút nhỏ hơn hoặc bằng"/"Rút hết" to match the target
# index layout (Có=142, Rút nhỏ hơn hoặc bằng=143, Rút hết=144).
$wsLoaiTietKiem.Range("G1").Value = "Có"
$wsLoaiTietKiem.Range("F1").Value = "Rút nhỏ hơn hoặc bằng"

$wsLoaiTietKiem.Range("F2").Value = "Rút hết"
$wsLoaiTietKiem.Range("G2").Value = "Có"

$wsLoaiTietKiem.Range("F3").Value = "Rút hết"
$wsLoaiTietKiem.Range("G3").Value = "Có"

# Row 1 grows to fit the longer text, column F widens a bit.
$wsLoaiTietKiem.Rows.Item(1).RowHeight = 45
$wsLoaiTietKiem.Columns.Item(6).ColumnWidth = 11.166666666666666

# ---------------------------------------------------------------------------
# SOTIETKIEM: column F becomes Text so large ID numbers stop getting mangled,
# and one bad value gets corrected.
# ---------------------------------------------------------------------------

# Fix the bad value first, *then* apply the Text format - doing it in the
# other order would make Excel store the literal as a shared string instead
# of a plain number (NumberFormat "@" coerces new values typed afterwards).
$wsSoTietKiem.Range("F3").Value = 15234634649
$wsSoTietKiem.Range("F1:F31").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping (must happen after the edits above so
# the final Select() calls are the ones that stick).
# ---------------------------------------------------------------------------

$wsLoaiTietKiem.Range("F2").Select() | Out-Null
$wsPhanQuyen.Range("A1:C18").Select() | Out-Null
$wsThamSo.Range("B3").Select() | Out-Null

# SOTIETKIEM becomes the active/visible tab, with G6 selected.
$wsSoTietKiem.Range("G6").Select() | Out-Null
